function Set-TextCell($ws, $ref, $val) {
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).NumberFormat = "General"
}

function Set-NumCell($ws, $ref, $val, $fmt) {
    if ($fmt) {
        $ws.Range($ref).NumberFormat = $fmt
    }
    $ws.Range($ref).Value = $val
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Weekly crime-stat table updates (rows 14-29) ---
Set-TextCell $ws "G14" "0"
Set-TextCell $ws "H14" "***.*"
Set-NumCell $ws "L14" -100 "#,##0.0;`"-`"#,##0.0"
Set-NumCell $ws "N14" -100 "#,##0.0;`"-`"#,##0.0"
Set-NumCell $ws "C15" 1 "#,##0"
Set-TextCell $ws "D15" "0"
Set-TextCell $ws "E15" "***.*"
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 300
$ws.Range("I15").Value = 3
$ws.Range("K15").Value = 200
$ws.Range("L15").Value = 200
Set-NumCell $ws "N15" 200 "#,##0.0;`"-`"#,##0.0"
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 34
$ws.Range("H16").Value = -23.529411764705
$ws.Range("I16").Value = 22
$ws.Range("J16").Value = 27
$ws.Range("K16").Value = -18.518518518518
$ws.Range("L16").Value = -15.384615384615
$ws.Range("M16").Value = 29.411764705882
$ws.Range("N16").Value = -74.418604651162
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -63.636363636363
$ws.Range("F17").Value = 30
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 40
$ws.Range("K17").Value = -37.5
$ws.Range("L17").Value = 4.166666666666
$ws.Range("M17").Value = 78.571428571428
$ws.Range("N17").Value = -51.923076923076
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 26
$ws.Range("H18").Value = -46.153846153846
$ws.Range("I18").Value = 10
$ws.Range("J18").Value = 22
$ws.Range("K18").Value = -54.545454545454
$ws.Range("L18").Value = -28.571428571428
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = -79.591836734693
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 100
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 28
$ws.Range("J19").Value = 19
$ws.Range("K19").Value = 47.368421052631
$ws.Range("L19").Value = 33.333333333333
$ws.Range("N19").Value = 0
$ws.Range("C20").Value = 10
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = 100
$ws.Range("F20").Value = 33
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = 73.684210526315
$ws.Range("I20").Value = 28
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 75
$ws.Range("L20").Value = 115.384615384615
$ws.Range("M20").Value = 211.111111111111
$ws.Range("N20").Value = -9.677419354838
$ws.Range("C21").Value = 35
$ws.Range("D21").Value = 41
$ws.Range("E21").Value = -14.634146341463
$ws.Range("F21").Value = 142
$ws.Range("G21").Value = 164
$ws.Range("H21").Value = -13.414634146341
$ws.Range("I21").Value = 116
$ws.Range("J21").Value = 125
$ws.Range("K21").Value = -7.2
$ws.Range("L21").Value = 16
$ws.Range("M21").Value = 93.333333333333
$ws.Range("N21").Value = -53.225806451612
Set-TextCell $ws "G22" "0"
Set-TextCell $ws "H22" "***.*"
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 9
$ws.Range("E23").Value = -66.666666666666
$ws.Range("G23").Value = 25
$ws.Range("H23").Value = -4
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 150
$ws.Range("M23").Value = 0
$ws.Range("C24").Value = 21
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 87
$ws.Range("H24").Value = 13.793103448275
$ws.Range("I24").Value = 64
$ws.Range("J24").Value = 62
$ws.Range("K24").Value = 3.225806451612
$ws.Range("L24").Value = 3.225806451612
$ws.Range("M24").Value = 23.076923076923
$ws.Range("C25").Value = 24
$ws.Range("D25").Value = 23
$ws.Range("E25").Value = 4.347826086956
$ws.Range("F25").Value = 76
$ws.Range("G25").Value = 69
$ws.Range("H25").Value = 10.144927536231
$ws.Range("I25").Value = 62
$ws.Range("J25").Value = 54
$ws.Range("K25").Value = 14.814814814814
$ws.Range("L25").Value = 34.782608695652
$ws.Range("M25").Value = 10.714285714285
Set-NumCell $ws "C26" 1 "#,##0"
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 4
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -20
$ws.Range("I26").Value = 3
$ws.Range("J26").Value = 3
$ws.Range("L26").Value = -25
$ws.Range("C27").Value = 3
Set-TextCell $ws "D27" "0"
Set-TextCell $ws "E27" "***.*"
$ws.Range("F27").Value = 9
$ws.Range("H27").Value = 80
$ws.Range("I27").Value = 7
$ws.Range("K27").Value = 133.333333333333
$ws.Range("L27").Value = 75
Set-TextCell $ws "C28" "0"
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 2
$ws.Range("L28").Value = -50
$ws.Range("M28").Value = -66.666666666666
$ws.Range("N28").Value = -75
Set-TextCell $ws "C29" "0"
$ws.Range("F29").Value = 2
$ws.Range("G29").Value = 2
$ws.Range("L29").Value = -50
$ws.Range("M29").Value = -66.666666666666
$ws.Range("N29").Value = -75
